# Updated symbol list on Mon Dec 12 10:42:02 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) refreshes plus the CEJI / BKEXToken
# row swap (rows 42 & 43) described by the scrape diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Price") cells get new quoted values straight from the feed.
# These look numeric, so Excel would normally coerce them to doubles on
# assignment (and in a couple of cases that rounds away a significant
# trailing zero, e.g. 0.08310 -> 0.0831). Marking the cells as Text first
# keeps them verbatim, matching how the sheet already stores every other
# column as text.
$priceUpdates = @{
    "D2"  = "282.79"
    "D3"  = "20.97"
    "D4"  = "6.226"
    "D5"  = "0.06178"
    "D6"  = "3.584"
    "D7"  = "6.558"
    "D8"  = "1.475"
    "D9"  = "0.8176"
    "D11" = "0.1643"
    "D12" = "0.08310"
    "D13" = "0.03614"
    "D15" = "0.09134"
    "D16" = "3.697"
    "D17" = "0.001638"
    "D18" = "0.04682"
    "D19" = "0.006435"
    "D20" = "0.006191"
    "D21" = "0.001067"
    "D23" = "3.819"
    "D40" = "0.04705"
    "D41" = "0.007104"
    "D44" = "0.01141"
    "D45" = "0.00006352"
    "D47" = "0.9987"
    "D48" = "0.002771"
    "D50" = "0.01238"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    # Restore the cell's original (default) formatting - only the text
    # content itself should change, not the look of the cell.
    $cell.ClearFormats()
}

# --- Rows 42 & 43: the ranking swapped BKEXToken and CEJI, each with a
# freshly scraped price (not simply each other's old price).
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1107"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003517"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "42CEJICEJI"
